# Daily attendance processing - 2025-10-20 19:17:15
# Reorder the "Recorded By" (column G) entries so that "System" (and other
# non-first authors) are listed first, matching the author order produced by
# the latest attendance sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "System, backup@backdoor.com, system"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    8   = "System, backup@backdoor.com"
    11  = "System, dnasr281@gmail.com"
    17  = "System, dnasr281@gmail.com"
    29  = "System, backup@backdoor.com, system"
    31  = "System, backup@backdoor.com"
    32  = "System, backup@backdoor.com"
    35  = "System, backup@backdoor.com"
    38  = "System, dnasr281@gmail.com"
    44  = "System, dnasr281@gmail.com"
    56  = "System, backup@backdoor.com, system"
    58  = "System, backup@backdoor.com"
    59  = "System, backup@backdoor.com"
    62  = "System, backup@backdoor.com"
    65  = "System, dnasr281@gmail.com"
    71  = "System, dnasr281@gmail.com"
    83  = "System, backup@backdoor.com"
    84  = "System, backup@backdoor.com"
    85  = "System, backup@backdoor.com"
    90  = "admin@admin.com, dnasr281@gmail.com"
    96  = "System, dnasr281@gmail.com"
    97  = "System, dnasr281@gmail.com"
    109 = "System, backup@backdoor.com"
    110 = "System, backup@backdoor.com"
    111 = "System, backup@backdoor.com"
    116 = "admin@admin.com, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    123 = "System, dnasr281@gmail.com"
    135 = "System, backup@backdoor.com"
    136 = "System, backup@backdoor.com"
    137 = "System, backup@backdoor.com"
    142 = "admin@admin.com, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
    149 = "System, dnasr281@gmail.com"
}

foreach ($rowNum in $updates.Keys) {
    $ws.Range("G$rowNum").Value = $updates[$rowNum]
}
